# This script reconstructs the full worksheet content, including the
# pre-existing rows (1-21) plus the newly added "generic" word_type column
# values and the new "stim details" block (rows 27-36) described by the diff.
#
# NOTE: the worksheet is rebuilt from a clean slate (Cells.Clear()) and all
# values (old and new) are written back in row-major / left-to-right order.
# This reproduces the exact final cell values of the target workbook; it is
# functionally equivalent to applying the diff in place (only four existing
# cells -- J2:J5 -- and the new rows 27:36 actually change content-wise, the
# rest of the OOXML diff is just shared-string table index churn caused by
# those insertions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Clear()

# Header row
$ws.Cells.Item(1,1).Value = 'number'
$ws.Cells.Item(1,2).Value = 'word'
$ws.Cells.Item(1,3).Value = 'kind'
$ws.Cells.Item(1,4).Value = 'carrier'
$ws.Cells.Item(1,5).Value = 'duplicate_image_filename'
$ws.Cells.Item(1,7).Value = 'order'
$ws.Cells.Item(1,8).Value = 'pair'
$ws.Cells.Item(1,9).Value = 'pair_words'
$ws.Cells.Item(1,10).Value = 'pair_kind'
$ws.Cells.Item(1,11).Value = 'carrier'

# Practice row p1
$ws.Cells.Item(2,1).Value = 'p1'
$ws.Cells.Item(2,3).Value = 'practice'
$ws.Cells.Item(2,7).Value = 7
$ws.Cells.Item(2,8).Value = 'A'
$ws.Cells.Item(2,9).Value = 'banana_kitty'
$ws.Cells.Item(2,10).Value = 'generic'
$ws.Cells.Item(2,11).Value = 'can'

# Practice row p2
$ws.Cells.Item(3,1).Value = 'p2'
$ws.Cells.Item(3,3).Value = 'practice'
$ws.Cells.Item(3,8).Value = 'B'
$ws.Cells.Item(3,9).Value = 'bear_cracker'
$ws.Cells.Item(3,10).Value = 'generic'
$ws.Cells.Item(3,11).Value = 'do'

# Practice row p3
$ws.Cells.Item(4,1).Value = 'p3'
$ws.Cells.Item(4,3).Value = 'practice'
$ws.Cells.Item(4,8).Value = 'C'
$ws.Cells.Item(4,9).Value = 'hair_cup'
$ws.Cells.Item(4,10).Value = 'generic'
$ws.Cells.Item(4,11).Value = 'look'

# Practice row p4
$ws.Cells.Item(5,1).Value = 'p4'
$ws.Cells.Item(5,3).Value = 'practice'
$ws.Cells.Item(5,8).Value = 'D'
$ws.Cells.Item(5,9).Value = 'cheerios_water'
$ws.Cells.Item(5,10).Value = 'generic'
$ws.Cells.Item(5,11).Value = 'where'

# Trial 1
$ws.Cells.Item(6,1).Value = 1
$ws.Cells.Item(6,2).Value = 'banana'
$ws.Cells.Item(6,3).Value = 'generic'
$ws.Cells.Item(6,4).Value = 'can'
$ws.Cells.Item(6,8).Value = 'E'
$ws.Cells.Item(6,11).Value = 'look'

# Trial 2
$ws.Cells.Item(7,1).Value = 2
$ws.Cells.Item(7,2).Value = 'kitty'
$ws.Cells.Item(7,3).Value = 'generic'
$ws.Cells.Item(7,4).Value = 'can'
$ws.Cells.Item(7,8).Value = 'F'
$ws.Cells.Item(7,11).Value = 'where'

# Trial 3
$ws.Cells.Item(8,1).Value = 3
$ws.Cells.Item(8,2).Value = 'bear'
$ws.Cells.Item(8,3).Value = 'generic'
$ws.Cells.Item(8,4).Value = 'do'
$ws.Cells.Item(8,8).Value = 'G'
$ws.Cells.Item(8,11).Value = 'can'

# Trial 4
$ws.Cells.Item(9,1).Value = 4
$ws.Cells.Item(9,2).Value = 'cracker'
$ws.Cells.Item(9,3).Value = 'generic'
$ws.Cells.Item(9,4).Value = 'do'
$ws.Cells.Item(9,8).Value = 'H'
$ws.Cells.Item(9,11).Value = 'do'

# Trial 5
$ws.Cells.Item(10,1).Value = 5
$ws.Cells.Item(10,2).Value = 'cup'
$ws.Cells.Item(10,3).Value = 'generic'
$ws.Cells.Item(10,4).Value = 'look'

# Trial 6
$ws.Cells.Item(11,1).Value = 6
$ws.Cells.Item(11,2).Value = 'hair'
$ws.Cells.Item(11,3).Value = 'generic'
$ws.Cells.Item(11,4).Value = 'look'

# Trial 7
$ws.Cells.Item(12,1).Value = 7
$ws.Cells.Item(12,2).Value = 'cheerios'
$ws.Cells.Item(12,3).Value = 'generic'
$ws.Cells.Item(12,4).Value = 'where'

# Trial 8
$ws.Cells.Item(13,1).Value = 8
$ws.Cells.Item(13,2).Value = 'water'
$ws.Cells.Item(13,3).Value = 'generic'
$ws.Cells.Item(13,4).Value = 'where'

$ws.Cells.Item(14,1).Value = 9

$ws.Cells.Item(15,1).Value = 10

$ws.Cells.Item(16,1).Value = 11

$ws.Cells.Item(17,1).Value = 12

$ws.Cells.Item(18,1).Value = 13

$ws.Cells.Item(19,1).Value = 14

$ws.Cells.Item(20,1).Value = 15

$ws.Cells.Item(21,1).Value = 16

# New block: stim details header
$ws.Cells.Item(27,1).Value = 'stim details'

# New block: stim details column headers
$ws.Cells.Item(28,1).Value = 'month'
$ws.Cells.Item(28,2).Value = 'word_type'
$ws.Cells.Item(28,3).Value = 'need_audio'
$ws.Cells.Item(28,4).Value = 'need_image'
$ws.Cells.Item(28,5).Value = 'word'
$ws.Cells.Item(28,6).Value = 'count'
$ws.Cells.Item(28,7).Value = 'find images'

# New block: stim details data
$ws.Cells.Item(29,1).Value = 6
$ws.Cells.Item(29,2).Value = 'video'

$ws.Cells.Item(30,1).Value = 6
$ws.Cells.Item(30,2).Value = 'video'

$ws.Cells.Item(31,1).Value = 7
$ws.Cells.Item(31,2).Value = 'video'

$ws.Cells.Item(32,1).Value = 7
$ws.Cells.Item(32,2).Value = 'video'

$ws.Cells.Item(33,1).Value = 6
$ws.Cells.Item(33,2).Value = 'audio'

$ws.Cells.Item(34,1).Value = 6
$ws.Cells.Item(34,2).Value = 'audio'

$ws.Cells.Item(35,1).Value = 7
$ws.Cells.Item(35,2).Value = 'audio'

$ws.Cells.Item(36,1).Value = 7
$ws.Cells.Item(36,2).Value = 'audio'
